$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new header "Space Column" in column E, row 1 (new shared string,
# extends the used range from A1:D6 to A1:E6).
$ws.Range("E1").Value = "Space Column"

# Reset the lingering stale selection (was "D15", outside the sheet's data
# range) back to the top-left cell.
[void]$ws.Range("A1").Select()
